$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update formulas for row 8 (C8 and D8 swap their multiplier flags)
$ws.Range("C8").Formula = "=`$K8*1.03*1.03*10%*0"
$ws.Range("D8").Formula = "=`$K8*1.03*1.03*1.03*10%*1"

# Update formula for row 10 (B10)
$ws.Range("B10").Formula = "=K9*(0.5/9)"

# Update the sheet view selection from C9 to D9
$ws.Activate()
$ws.Range("D9").Select()
